$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data (row 3). The Artwork Url / Lid Url columns reuse
# the same values as row 2, so copy those two cells across (this preserves
# the shared-string reuse and the special font formatting applied to the
# Lid Url column).
$ws.Range("A3").Value = "US-Social Beverages"
$ws.Range("B3").Value = "Primary"
$ws.Range("C2:D2").Copy($ws.Range("C3:D3"))

# Update the selection to match the new active cell/range
$ws.Range("B3:D3").Select() | Out-Null

# Column A widens (best-fit) to accommodate the longer new text
$ws.Columns.Item(1).ColumnWidth = 18.14
